$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The Price column (D) stores numeric-looking values as plain text in the
# source workbook (e.g. "0.9978", "5.021"). Mark each such cell as Text
# before writing so Excel does not silently convert it to a number and
# drop significant trailing zeros / switch to scientific notation.

$ws.Range("D2").Value = '29.335.85'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '1.846.48'
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9978'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.98'
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6268'
$ws.Range("E6").Value = '  -0.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9981'
$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07606'
$ws.Range("E8").Value = '  -1.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2902'
$ws.Range("E9").Value = '  -1.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.71'
$ws.Range("E10").Value = '  +0.83%  '

$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.021'
$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6781'
$ws.Range("E13").Value = '  -0.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001049'
$ws.Range("E14").Value = '  -3.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.99'
$ws.Range("E15").Value = '  -0.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.132'
$ws.Range("E16").Value = '  -0.33%  '

$ws.Range("D17").Value = '29.345.03'
$ws.Range("E17").Value = '  -0.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.48'
$ws.Range("E18").Value = '  -0.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.32'
$ws.Range("E19").Value = '  -1.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9980'
$ws.Range("E20").Value = '  -0.24%  '

$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9988'
$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '158.03'
$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1383'
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.402'
$ws.Range("E25").Value = '  +0.45%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.64'
$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.402'
$ws.Range("E27").Value = '  +6.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.458'
$ws.Range("E28").Value = '  -0.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05599'
$ws.Range("E29").Value = '  -1.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.104'
$ws.Range("E30").Value = '  -0.17%  '

$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.161'
$ws.Range("E32").Value = '  +0.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.832'
$ws.Range("E33").Value = '  -0.97%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.6963'

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.585'
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01801'
$ws.Range("E36").Value = '  +0.13%  '

$ws.Range("D37").Value = '1.225.88'
$ws.Range("E37").Value = '  -0.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.719'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.350'
$ws.Range("E39").Value = '  -1.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9029'
$ws.Range("E40").Value = '  -1.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9981'
$ws.Range("E41").Value = '  -0.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.17'
$ws.Range("E42").Value = '  -0.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.53'
$ws.Range("E43").Value = '  -1.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.194'
$ws.Range("E44").Value = '  +0.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3990'
$ws.Range("E45").Value = '  -0.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '9.020'
$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.674'
$ws.Range("E47").Value = '  -0.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.1138'
$ws.Range("E48").Value = '  +1.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.05701'
$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.4619'
$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.338'
$ws.Range("E51").Value = '  -1.18%  '
